$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.956.30"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.383.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.63"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.38%  "

$ws.Range("E10").Value = "  -1.01%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.963.12"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("E13").Value = "  +2.17%  "

$ws.Range("E14").Value = "  -2.03%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.383.01"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.072.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.08"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.44"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.551"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.41%  "

$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("E25").Value = "  -4.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.522.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.28"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.94"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.414.19"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.65%  "

$ws.Range("E39").Value = "  -4.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0766"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.02"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.38"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("E45").Value = "  -2.89%  "

$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.455.87"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.99"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.71"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.46%  "

$ws.Range("E50").Value = "  +2.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.13"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.81%  "
